$wb = $excel.ActiveWorkbook

# The edit is made on the "Repayment schedule" sheet (already the active /
# selected tab in this workbook).
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column in front of column N (the "Late" column),
# pushing "Late" / the blank-heading column / "Outstanding" one column to
# the right (N->O, O->P, P->Q). This is a plain column insert, so the new
# column inherits the formatting (and width) of the column immediately to
# its left (column M).
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Leave the selection where the author left it after making the edit.
$ws.Range("J16").Select()
